# Modify spacing before/after for Heading2 and Heading3 styles in the
# document's style sheet: before=200/after=0 (twips) -> before=240/after=240
# (twips). Word COM expresses ParagraphFormat.SpaceBefore/SpaceAfter in
# points, and 1 point == 20 twips, so twips/20 gives the point value to set.
$d = $word.ActiveDocument

$heading2 = $d.Styles("Heading2")
$heading2.ParagraphFormat.SpaceBefore = 12
$heading2.ParagraphFormat.SpaceAfter = 12

$heading3 = $d.Styles("Heading3")
$heading3.ParagraphFormat.SpaceBefore = 12
$heading3.ParagraphFormat.SpaceAfter = 12

Write-Output "Updated Heading2/Heading3 spacing (before=240 twips, after=240 twips)"
